# feat: add 2022-Q1 data
#
# The workbook has 4 sheets: 2021-Q2, 2021-Q3, 2021-Q4, 总计.
# We add a new "2022-Q1" detail sheet (taking over the slot/sheetId that
# used to belong to "总计") and a fresh "总计" sheet at the end that keeps
# the historical summary rows plus a new one for 2022-Q1.

$wb = $excel.ActiveWorkbook

$totalOrig = $wb.Worksheets.Item("总计")

# Duplicate the existing "总计" sheet so we end up with two sheets holding
# the old summary content; one keeps the name "总计" (becomes the new
# sheet, at the end) and the other is turned into the new "2022-Q1" detail
# sheet (reusing the original sheet's slot, right where "总计" used to be).
$totalOrig.Copy($null, $totalOrig)

$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

$newTotal = $wb.Worksheets.Item("总计 (2)")
$newTotal.Name = "总计"

# ---------------------------------------------------------------------
# Build the "2022-Q1" fund-holdings detail sheet
# ---------------------------------------------------------------------
$q1.Cells.Clear()

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Copy the header cell style (bold, centered, bordered - style index 2 in
# the original workbook) from sheet "2021-Q4" onto the new header row, and
# also copy the index-column style onto column A for the data rows.
$styleSrc = $wb.Worksheets.Item("2021-Q4")
$styleSrc.Range("B1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$indexStyleSrc = $styleSrc.Range("A2")

$fundRows = @(
  @{code="001410"; name="信达澳银新能源产业股票";             scale="140.41"; pos="92.06"; ratio="1.54"; mv="2.1623"; rank=6},
  @{code="012608"; name="信达澳银领先智选混合型证券投资基金"; scale="38.78";  pos="90.57"; ratio="1.78"; mv="0.6903"; rank=6},
  @{code="011188"; name="信达澳银星奕混合A";                  scale="35.73";  pos="92.14"; ratio="1.78"; mv="0.6360"; rank=6},
  @{code="006257"; name="信达澳银先进智造股票";               scale="21.53";  pos="93.84"; ratio="1.69"; mv="0.3639"; rank=6},
  @{code="011223"; name="信达澳银星奕混合C";                  scale="9.21";   pos="92.14"; ratio="1.78"; mv="0.1639"; rank=6},
  @{code="009511"; name="信达澳银研究优选混合";               scale="9.41";   pos="92.12"; ratio="1.72"; mv="0.1619"; rank=6},
  @{code="009437"; name="信达澳银科技创新一年定期开放混合A"; scale="5.57";   pos="94.43"; ratio="1.78"; mv="0.0991"; rank=6},
  @{code="005632"; name="鹏华量化先锋混合";                   scale="3.10";   pos="92.91"; ratio="1.71"; mv="0.0530"; rank=5},
  @{code="014133"; name="工银中证500六个月持有指数增强A";     scale="3.07";   pos="93.69"; ratio="1.18"; mv="0.0362"; rank=1},
  @{code="014344"; name="鹏华中证500指数增强A";               scale="1.99";   pos="92.63"; ratio="1.65"; mv="0.0328"; rank=9},
  @{code="009438"; name="信达澳银科技创新一年定期开放混合C"; scale="0.84";   pos="94.43"; ratio="1.78"; mv="0.0150"; rank=6},
  @{code="014134"; name="工银中证500六个月持有指数增强C";     scale="1.12";   pos="93.69"; ratio="1.18"; mv="0.0132"; rank=1},
  @{code="014345"; name="鹏华中证500指数增强C";               scale="0.78";   pos="92.63"; ratio="1.65"; mv="0.0129"; rank=9}
)

$r = 2
foreach ($fr in $fundRows) {
  $q1.Range("A$r").Value = $r - 2
  $indexStyleSrc.Copy()
  $q1.Range("A$r").PasteSpecial(-4122)

  # Fund code / name / scale / position / ratio / market-value columns are
  # stored as text in the source data (e.g. "140.41"), so force the
  # number format to Text before writing, otherwise Excel would silently
  # reinterpret them as numbers (and fund codes would lose leading zeros).
  $q1.Range("B${r}:G${r}").NumberFormat = "@"

  $q1.Range("B$r").Value = $fr.code
  $q1.Range("C$r").Value = $fr.name
  $q1.Range("D$r").Value = $fr.scale
  $q1.Range("E$r").Value = $fr.pos
  $q1.Range("F$r").Value = $fr.ratio
  $q1.Range("G$r").Value = $fr.mv
  $q1.Range("H$r").Value = $fr.rank

  $r++
}

# ---------------------------------------------------------------------
# Build the new "总计" summary sheet: a new 2022-Q1 row on top of the
# historical rows.
# ---------------------------------------------------------------------
$newTotal.Cells.Clear()

$newTotal.Range("B1").Value = "日期"
$newTotal.Range("C1").Value = "持有数量(只)"
$newTotal.Range("D1").Value = "持有市值(亿元)"
$styleSrc.Range("B1").Copy()
$newTotal.Range("B1:D1").PasteSpecial(-4122)

$summaryRows = @(
  @{date="2022-Q1"; count=13; mv=4.44},
  @{date="2021-Q4"; count=11; mv=5.63},
  @{date="2021-Q3"; count=11; mv=4.3},
  @{date="2021-Q2"; count=1;  mv=0.03}
)

$r = 2
foreach ($sr in $summaryRows) {
  $newTotal.Range("A$r").Value = $r - 2
  $indexStyleSrc.Copy()
  $newTotal.Range("A$r").PasteSpecial(-4122)

  $newTotal.Range("B$r").Value = $sr.date
  $newTotal.Range("C$r").Value = $sr.count
  $newTotal.Range("D$r").Value = $sr.mv

  $r++
}
